$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 385
$ws.Range("E2").Value = 42
$ws.Range("F2").Value = 42
$ws.Range("G2").Value = 42
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 30
$ws.Range("K2").Value = 394
$ws.Range("L2").Value = 62
$ws.Range("M2").Value = 332
$ws.Range("N2").Value = 332
$ws.Range("P2").Value = 65
$ws.Range("Q2").Value = 33
$ws.Range("R2").Value = -2
$ws.Range("S2").Value = -21
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = 26
$ws.Range("V2").Value = 9
$ws.Range("W2").Value = 10.79
$ws.Range("X2").Value = 7.91
$ws.Range("Y2").Value = 9.300000000000001
$ws.Range("Z2").Value = 7.78
$ws.Range("AA2").Value = 18.72
$ws.Range("AB2").Value = 430.32
$ws.Range("AC2").Value = 234
$ws.Range("AD2").Value = 10.95
$ws.Range("AE2").Value = 2555
$ws.Range("AF2").Value = 1
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 5.85
$ws.Range("AI2").Value = 64.06
$ws.Range("AJ2").Value = 13000000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 437
$ws.Range("E3").Value = 70
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = 70
$ws.Range("H3").Value = 53
$ws.Range("I3").Value = 53
$ws.Range("K3").Value = 431
$ws.Range("L3").Value = 65
$ws.Range("M3").Value = 365
$ws.Range("N3").Value = 365
$ws.Range("P3").Value = 65
$ws.Range("Q3").Value = 57
$ws.Range("R3").Value = -24
$ws.Range("S3").Value = -22
$ws.Range("T3").Value = 15
$ws.Range("U3").Value = 42
$ws.Range("V3").Value = 6
$ws.Range("W3").Value = 15.94
$ws.Range("X3").Value = 12.15
$ws.Range("Y3").Value = 15.23
$ws.Range("Z3").Value = 12.88
$ws.Range("AA3").Value = 17.88
$ws.Range("AB3").Value = 480.37
$ws.Range("AC3").Value = 409
$ws.Range("AD3").Value = 7.39
$ws.Range("AE3").Value = 2809
$ws.Range("AF3").Value = 1.07
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 6.62
$ws.Range("AI3").Value = 48.94
$ws.Range("AJ3").Value = 13000000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 424
$ws.Range("E4").Value = 52
$ws.Range("F4").Value = 52
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 13
$ws.Range("K4").Value = 466
$ws.Range("L4").Value = 112
$ws.Range("M4").Value = 354
$ws.Range("N4").Value = 354
$ws.Range("P4").Value = 65
$ws.Range("Q4").Value = 48
$ws.Range("R4").Value = -32
$ws.Range("S4").Value = -26
$ws.Range("T4").Value = 43
$ws.Range("U4").Value = 5
$ws.Range("V4").Value = 8
$ws.Range("W4").Value = 12.2
$ws.Range("X4").Value = 3.16
$ws.Range("Y4").Value = 3.73
$ws.Range("Z4").Value = 2.99
$ws.Range("AA4").Value = 31.54
$ws.Range("AB4").Value = 462.14
$ws.Range("AC4").Value = 103
$ws.Range("AD4").Value = 42.9
$ws.Range("AE4").Value = 2726
$ws.Range("AF4").Value = 1.63
$ws.Range("AG4").Value = 175
$ws.Range("AH4").Value = 3.95
$ws.Range("AI4").Value = 169.48
$ws.Range("AJ4").Value = 13000000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 448
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 18
$ws.Range("H5").Value = 17
$ws.Range("I5").Value = 17
$ws.Range("K5").Value = 469
$ws.Range("L5").Value = 127
$ws.Range("M5").Value = 342
$ws.Range("N5").Value = 342
$ws.Range("P5").Value = 65
$ws.Range("Q5").Value = 41
$ws.Range("R5").Value = -40
$ws.Range("S5").Value = -5
$ws.Range("T5").Value = 45
$ws.Range("U5").Value = -4
$ws.Range("V5").Value = 26
$ws.Range("W5").Value = 4.12
$ws.Range("X5").Value = 3.71
$ws.Range("Y5").Value = 4.78
$ws.Range("Z5").Value = 3.56
$ws.Range("AA5").Value = 37.04
$ws.Range("AB5").Value = 450.67
$ws.Range("AC5").Value = 128
$ws.Range("AD5").Value = 22.81
$ws.Range("AE5").Value = 2632
$ws.Range("AF5").Value = 1.11
$ws.Range("AG5").Value = 125
$ws.Range("AH5").Value = 4.28
$ws.Range("AI5").Value = 97.64
$ws.Range("AJ5").Value = 13000000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 594
$ws.Range("E6").Value = 30
$ws.Range("F6").Value = 30
$ws.Range("G6").Value = 29
$ws.Range("H6").Value = 24
$ws.Range("I6").Value = 24
$ws.Range("K6").Value = 507
$ws.Range("L6").Value = 153
$ws.Range("M6").Value = 354
$ws.Range("N6").Value = 354
$ws.Range("P6").Value = 65
$ws.Range("Q6").Value = 40
$ws.Range("R6").Value = -25
$ws.Range("S6").Value = -8
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = 35
$ws.Range("V6").Value = 31
$ws.Range("W6").Value = 5.09
$ws.Range("X6").Value = 4.06
$ws.Range("Y6").Value = 6.93
$ws.Range("Z6").Value = 4.94
$ws.Range("AA6").Value = 43.14
$ws.Range("AB6").Value = 468.05
$ws.Range("AC6").Value = 185
$ws.Range("AD6").Value = 28.25
$ws.Range("AE6").Value = 2724
$ws.Range("AF6").Value = 1.92
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 1.91
$ws.Range("AI6").Value = 53.91
$ws.Range("AJ6").Value = 13000000

# Row 7: remove all financial data, keep A/B/C
$ws.Range("D7:AJ7").ClearContents()

# Row 8: remove all financial data, keep A/B/C
$ws.Range("D8:AJ8").ClearContents()

# Row 9: remove all financial data, keep A/B/C
$ws.Range("D9:AJ9").ClearContents()
